$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The source edit swaps the data of row 20 and row 21 (two observation
# records traded places). Several columns happen to hold equal values
# in both rows (D, I, T, U, V, W, Y, AA, AD, AE, AG, AT, AY) so only the
# cells below actually change value and need to be written.
# ---------------------------------------------------------------------

# --- Capture the current (pre-edit) values of both rows ---
$A20 = $ws.Range("A20").Value()
$B20 = $ws.Range("B20").Value()
$E20 = $ws.Range("E20").Value()
$F20 = $ws.Range("F20").Value()
$G20 = $ws.Range("G20").Value()
$H20 = $ws.Range("H20").Value()
$P20 = $ws.Range("P20").Value()
$Q20 = $ws.Range("Q20").Value()
$R20 = $ws.Range("R20").Value()
$S20 = $ws.Range("S20").Value()
$Z20 = $ws.Range("Z20").Value()
$AB20 = $ws.Range("AB20").Value()
$AW20 = $ws.Range("AW20").Value()
$AX20 = $ws.Range("AX20").Value()

$A21 = $ws.Range("A21").Value()
$B21 = $ws.Range("B21").Value()
$E21 = $ws.Range("E21").Value()
$F21 = $ws.Range("F21").Value()
$G21 = $ws.Range("G21").Value()
$H21 = $ws.Range("H21").Value()
$M21 = $ws.Range("M21").Value()
$P21 = $ws.Range("P21").Value()
$Q21 = $ws.Range("Q21").Value()
$R21 = $ws.Range("R21").Value()
$S21 = $ws.Range("S21").Value()
$Z21 = $ws.Range("Z21").Value()
$AB21 = $ws.Range("AB21").Value()
$AC21 = $ws.Range("AC21").Value()
$AW21 = $ws.Range("AW21").Value()
$AX21 = $ws.Range("AX21").Value()

# --- Row 20 becomes what row 21 used to hold ---
$ws.Range("A20").Value = $A21
$ws.Range("B20").Value = $B21
$ws.Range("E20").Value = $E21
$ws.Range("F20").Value = $F21
$ws.Range("G20").Value = $G21
$ws.Range("H20").Value = $H21
$ws.Range("M20").Value = $M21
$ws.Range("P20").Value = $P21
$ws.Range("Q20").Value = $Q21
$ws.Range("R20").Value = $R21
$ws.Range("S20").Value = $S21
$ws.Range("Z20").Value = $Z21
$ws.Range("AB20").Value = $AB21
$ws.Range("AC20").Value = $AC21
$ws.Range("AW20").Value = $AW21
$ws.Range("AX20").Value = $AX21

# --- Row 21 becomes what row 20 used to hold ---
$ws.Range("A21").Value = $A20
$ws.Range("B21").Value = $B20
$ws.Range("E21").Value = $E20
$ws.Range("F21").Value = $F20
$ws.Range("G21").Value = $G20
$ws.Range("H21").Value = $H20
$ws.Range("M21").ClearContents()
$ws.Range("P21").Value = $P20
$ws.Range("Q21").Value = $Q20
$ws.Range("R21").Value = $R20
$ws.Range("S21").Value = $S20
$ws.Range("Z21").Value = $Z20
$ws.Range("AB21").Value = $AB20
$ws.Range("AC21").ClearContents()
$ws.Range("AW21").Value = $AW20
$ws.Range("AX21").Value = $AX20
